$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 2001
$ws.Range("I2").Value = 2001
$ws.Range("K2").Value = 2001
$ws.Range("M2").Value = -1888
# Row 11
$ws.Range("H11").Value = 376.625
$ws.Range("I11").Value = 376.625
$ws.Range("K11").Value = 376.625
$ws.Range("M11").Value = -236.625
# Row 12
$ws.Range("H12").Value = 868.875
$ws.Range("J12").Value = 883.6667
$ws.Range("L12").Value = 883.6667
$ws.Range("N12").Value = -1223.6667
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
# Row 138
$ws.Range("H138").Value = 6189.143
$ws.Range("J138").Value = 6637.952
$ws.Range("L138").Value = 19913.856
$ws.Range("N138").Value = -30193.856

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2331.9092
$ws.Range("I2").Value = 2331.9092
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2331.9092
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2218.9092
$ws.Range("N2").ClearContents()
# Row 110
$ws.Range("H110").Value = 250000000
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
# Row 116
$ws.Range("H116").Value = 2331.9092
$ws.Range("I116").Value = 2331.9092
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2331.9092
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -37.90920000000006
$ws.Range("N116").ClearContents()
# Row 122
$ws.Range("H122").Value = 6207.048
$ws.Range("J122").Value = 9699.9
$ws.Range("L122").Value = 29099.7
$ws.Range("N122").Value = -33999.7

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2331.9092
$ws.Range("I3").Value = 2331.9092
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2331.9092
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2217.9092
$ws.Range("N3").ClearContents()
# Row 86
$ws.Range("H86").Value = 3099.75
$ws.Range("I86").Value = 3099.75
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3099.75
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1976.75
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 3099.75
$ws.Range("I89").Value = 3099.75
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 15498.75
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -9882.75
$ws.Range("N89").ClearContents()
# Row 94
$ws.Range("H94").Value = 368752.66
$ws.Range("I94").Value = 551254
$ws.Range("K94").Value = 551254
$ws.Range("M94").Value = -550803
# Row 99
$ws.Range("H99").Value = 2346.182
$ws.Range("I99").Value = 2346.182
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2346.182
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -848.1819999999998
$ws.Range("N99").ClearContents()
# Row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1596.75
$ws.Range("I16").Value = 1495.6666
$ws.Range("J16").Value = 1900
$ws.Range("K16").Value = 1495.6666
$ws.Range("L16").Value = 1900
$ws.Range("M16").Value = -1208.6666
$ws.Range("N16").Value = -2474
# Row 74
$ws.Range("H74").Value = 41749.668
$ws.Range("J74").Value = 41874.5
$ws.Range("L74").Value = 41874.5
$ws.Range("N74").Value = -43622.5
# Row 77
$ws.Range("H77").Value = 41749.668
$ws.Range("J77").Value = 41874.5
$ws.Range("L77").Value = 125623.5
$ws.Range("N77").Value = -134359.5
# Row 107
$ws.Range("H107").Value = 713.4
$ws.Range("I107").Value = 463.9091
$ws.Range("J107").Value = 1399.5
$ws.Range("K107").Value = 463.9091
$ws.Range("L107").Value = 1399.5
$ws.Range("M107").Value = 1456.0909
$ws.Range("N107").Value = -5239.5
# Row 113
$ws.Range("H113").Value = 1596.75
$ws.Range("I113").Value = 1495.6666
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 1495.6666
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 674.3334
$ws.Range("N113").Value = -6240
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 141
$ws.Range("H141").Value = 577264.8
$ws.Range("J141").Value = 577264.8
$ws.Range("L141").Value = 577264.8
$ws.Range("N141").Value = -587624.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 994.1429000000001
$ws.Range("I131").Value = 994.1429000000001
$ws.Range("K131").Value = 2982.4287
$ws.Range("M131").Value = 2057.5713

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3700
$ws.Range("I80").Value = 3625
$ws.Range("K80").Value = 3625
$ws.Range("M80").Value = -2627
# Row 83
$ws.Range("H83").Value = 3700
$ws.Range("I83").Value = 3625
$ws.Range("K83").Value = 18125
$ws.Range("M83").Value = -13133
# Row 107
$ws.Range("H107").Value = 37038416
$ws.Range("J107").Value = 66668396
$ws.Range("L107").Value = 66668396
$ws.Range("N107").Value = -66672236
# Row 113
$ws.Range("H113").Value = 3900
$ws.Range("I113").Value = 3900
$ws.Range("K113").Value = 3900
$ws.Range("M113").Value = -1730
# Row 122
$ws.Range("H122").Value = 1274.1428
$ws.Range("I122").Value = 1274.1428
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3822.4284
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1372.4284
$ws.Range("N122").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 4000
$ws.Range("I22").Value = 5875
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 5875
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = -5580
$ws.Range("N22").Value = -840
# Row 27
$ws.Range("H27").Value = 4000
$ws.Range("I27").Value = 5875
$ws.Range("J27").Value = 250
$ws.Range("K27").Value = 5875
$ws.Range("L27").Value = 250
$ws.Range("M27").Value = -5768
$ws.Range("N27").Value = -464
# Row 135
$ws.Range("H135").Value = 199999
$ws.Range("J135").Value = 199999
$ws.Range("L135").Value = 199999
$ws.Range("N135").Value = -210139

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 4460.2
$ws.Range("I126").Value = 2099.5715
$ws.Range("K126").Value = 6298.7145
$ws.Range("M126").Value = -3828.7145
